$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1964.1471
$ws.Range("I62").Value = 1923
$ws.Range("J62").Value = 2078.4443
$ws.Range("K62").Value = 1923
$ws.Range("L62").Value = 2078.4443
$ws.Range("M62").Value = -1299
$ws.Range("N62").Value = -3326.4443
$ws.Range("H65").Value = 1964.1471
$ws.Range("I65").Value = 1923
$ws.Range("J65").Value = 2078.4443
$ws.Range("K65").Value = 9615
$ws.Range("L65").Value = 10392.2215
$ws.Range("M65").Value = -6495
$ws.Range("N65").Value = -16632.2215
$ws.Range("H70").Value = 1683.1428
$ws.Range("I70").Value = 1296.3846
$ws.Range("J70").Value = 2018.3334
$ws.Range("K70").Value = 3889.1538
$ws.Range("L70").Value = 6055.0002
$ws.Range("M70").Value = -3619.1538
$ws.Range("N70").Value = -6595.0002
$ws.Range("H73").Value = 1683.1428
$ws.Range("I73").Value = 1296.3846
$ws.Range("J73").Value = 2018.3334
$ws.Range("K73").Value = 3889.1538
$ws.Range("L73").Value = 6055.0002
$ws.Range("M73").Value = -2953.1538
$ws.Range("N73").Value = -7927.0002
$ws.Range("H75").Value = 23535
$ws.Range("J75").Value = 23535
$ws.Range("L75").Value = 23535
$ws.Range("N75").Value = -25407
$ws.Range("H78").Value = 23535
$ws.Range("J78").Value = 23535
$ws.Range("L78").Value = 70605
$ws.Range("N78").Value = -79965
$ws.Range("H94").Value = 2143.75
$ws.Range("I94").Value = 2143.75
$ws.Range("K94").Value = 2143.75
$ws.Range("M94").Value = -1692.75
$ws.Range("H129").Value = 1037.7609
$ws.Range("I129").Value = 389.4
$ws.Range("J129").Value = 1116.8292
$ws.Range("K129").Value = 1168.2
$ws.Range("L129").Value = 3350.487599999999
$ws.Range("M129").Value = 3831.8
$ws.Range("N129").Value = -13350.4876
$ws.Range("H137").Value = 954.72
$ws.Range("I137").Value = 832.5263
$ws.Range("J137").Value = 1341.6666
$ws.Range("K137").Value = 2497.5789
$ws.Range("L137").Value = 4024.9998
$ws.Range("M137").Value = 52.42110000000002
$ws.Range("N137").Value = -9124.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8333.166999999999
$ws.Range("J63").Value = 8400
$ws.Range("L63").Value = 8400
$ws.Range("N63").Value = -9772
$ws.Range("H66").Value = 8333.166999999999
$ws.Range("J66").Value = 8400
$ws.Range("L66").Value = 42000
$ws.Range("N66").Value = -48864
$ws.Range("H74").Value = 1921.2
$ws.Range("I74").Value = 1702.4
$ws.Range("J74").Value = 2140
$ws.Range("K74").Value = 1702.4
$ws.Range("L74").Value = 2140
$ws.Range("M74").Value = -828.4000000000001
$ws.Range("N74").Value = -3888
$ws.Range("H77").Value = 1921.2
$ws.Range("I77").Value = 1702.4
$ws.Range("J77").Value = 2140
$ws.Range("K77").Value = 8512
$ws.Range("L77").Value = 10700
$ws.Range("M77").Value = -4144
$ws.Range("N77").Value = -19436
$ws.Range("H128").Value = 34999.5
$ws.Range("J128").Value = 34999.5
$ws.Range("L128").Value = 34999.5
$ws.Range("N128").Value = -44959.5
$ws.Range("H132").Value = 47668284
$ws.Range("I132").Value = 111111860
$ws.Range("J132").Value = 85603.5
$ws.Range("K132").Value = 333335580
$ws.Range("L132").Value = 256810.5
$ws.Range("M132").Value = -333333050
$ws.Range("N132").Value = -261870.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1179
$ws.Range("I36").Value = 1018.5
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 1018.5
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -484.5
$ws.Range("N36").Value = -2568
$ws.Range("H42").Value = 55336.363
$ws.Range("J42").Value = 55336.363
$ws.Range("L42").Value = 55336.363
$ws.Range("N42").Value = -55992.363
$ws.Range("H82").Value = 11993.375
$ws.Range("I82").Value = 2221.3333
$ws.Range("J82").Value = 17856.6
$ws.Range("K82").Value = 2221.3333
$ws.Range("L82").Value = 17856.6
$ws.Range("M82").Value = -1838.3333
$ws.Range("N82").Value = -18622.6
$ws.Range("H85").Value = 11993.375
$ws.Range("I85").Value = 2221.3333
$ws.Range("J85").Value = 17856.6
$ws.Range("K85").Value = 2221.3333
$ws.Range("L85").Value = 17856.6
$ws.Range("M85").Value = -895.3332999999998
$ws.Range("N85").Value = -20508.6
$ws.Range("H86").Value = 30335608
$ws.Range("J86").Value = 113480.89
$ws.Range("L86").Value = 113480.89
$ws.Range("N86").Value = -115726.89
$ws.Range("H89").Value = 30335608
$ws.Range("J89").Value = 113480.89
$ws.Range("L89").Value = 567404.45
$ws.Range("N89").Value = -578636.45
$ws.Range("H134").Value = 9416.541999999999
$ws.Range("I134").Value = 3335.6875
$ws.Range("J134").Value = 21578.25
$ws.Range("K134").Value = 10007.0625
$ws.Range("L134").Value = 64734.75
$ws.Range("M134").Value = -7472.0625
$ws.Range("N134").Value = -69804.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7145137
$ws.Range("I31").Value = 9261011
$ws.Range("J31").Value = 4062.5
$ws.Range("K31").Value = 9261011
$ws.Range("L31").Value = 4062.5
$ws.Range("M31").Value = -9260716
$ws.Range("N31").Value = -4652.5
$ws.Range("H34").Value = 7145137
$ws.Range("I34").Value = 9261011
$ws.Range("J34").Value = 4062.5
$ws.Range("K34").Value = 9261011
$ws.Range("L34").Value = 4062.5
$ws.Range("M34").Value = -9260809
$ws.Range("N34").Value = -4466.5
$ws.Range("H58").Value = 1279.037
$ws.Range("I58").Value = 1213
$ws.Range("J58").Value = 1411.1111
$ws.Range("K58").Value = 1213
$ws.Range("L58").Value = 1411.1111
$ws.Range("M58").Value = -1010
$ws.Range("N58").Value = -1817.1111
$ws.Range("H99").Value = 263025.78
$ws.Range("I99").Value = 348340.25
$ws.Range("K99").Value = 348340.25
$ws.Range("M99").Value = -346842.25
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -54868
$ws.Range("H126").Value = 263025.78
$ws.Range("I126").Value = 348340.25
$ws.Range("K126").Value = 1045020.75
$ws.Range("M126").Value = -1042550.75
$ws.Range("H136").Value = 1279.037
$ws.Range("I136").Value = 1213
$ws.Range("J136").Value = 1411.1111
$ws.Range("K136").Value = 3639
$ws.Range("L136").Value = 4233.3333
$ws.Range("M136").Value = -1089
$ws.Range("N136").Value = -9333.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2783.3333
$ws.Range("J106").Value = 2783.3333
$ws.Range("L106").Value = 8349.999899999999
$ws.Range("N106").Value = -10241.9999
$ws.Range("H117").Value = 102003.2
$ws.Range("J117").Value = 145276
$ws.Range("L117").Value = 435828
$ws.Range("N117").Value = -442712
$ws.Range("H129").Value = 13335434
$ws.Range("I129").Value = 4766.6665
$ws.Range("J129").Value = 15153253
$ws.Range("K129").Value = 14299.9995
$ws.Range("L129").Value = 45459759
$ws.Range("M129").Value = -9299.999500000002
$ws.Range("N129").Value = -45469759

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 424898
$ws.Range("I132").Value = 56902
$ws.Range("K132").Value = 170706
$ws.Range("M132").Value = -168176
$ws.Range("H139").Value = 25000
$ws.Range("J139").Value = 25000
$ws.Range("L139").Value = 25000
$ws.Range("N139").Value = -35280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 5300
$ws.Range("J42").Value = 5300
$ws.Range("L42").Value = 5300
$ws.Range("N42").Value = -6426
$ws.Range("H44").Value = 6000
$ws.Range("J44").Value = 6000
$ws.Range("L44").Value = 6000
$ws.Range("N44").Value = -7108
$ws.Range("H49").Value = 5300
$ws.Range("J49").Value = 5300
$ws.Range("L49").Value = 5300
$ws.Range("N49").Value = -5594
$ws.Range("H80").Value = 28000
$ws.Range("J80").Value = 28000
$ws.Range("L80").Value = 28000
$ws.Range("N80").Value = -30246
$ws.Range("H83").Value = 28000
$ws.Range("J83").Value = 28000
$ws.Range("L83").Value = 84000
$ws.Range("N83").Value = -95232
$ws.Range("H100").Value = 2124.8147
$ws.Range("I100").Value = 1619.0526
$ws.Range("J100").Value = 3326
$ws.Range("K100").Value = 1619.0526
$ws.Range("L100").Value = 3326
$ws.Range("M100").Value = -1078.0526
$ws.Range("N100").Value = -4408

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 92633.91
$ws.Range("I96").Value = 1281.8572
$ws.Range("J96").Value = 252500
$ws.Range("K96").Value = 1281.8572
$ws.Range("L96").Value = 252500
$ws.Range("M96").Value = 91.14280000000008
$ws.Range("N96").Value = -255246
$ws.Range("H126").Value = 100002120
$ws.Range("I126").Value = 142858750
$ws.Range("K126").Value = 428576250
$ws.Range("M126").Value = -428573780
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
